# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rewrites the worker mora table (B16:G27 on "Hoja1") with the updated /
# re-sorted dataset: a new worker (WILLIAM ROJAS SANCHEZ, 79329212) is
# placed first with period 1702, and the remaining workers are re-ordered
# so period 1703 rows come before period 1704 rows, with several Valor
# Mora (F) / Salario Basico (G) amounts updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $tipoDoc, $numDoc, $nombre, $periodo, $valorMora, $salario) {
    $ws.Range("B$row").Value = $tipoDoc
    $ws.Range("C$row").Value = $numDoc
    $ws.Range("D$row").Value = $nombre
    $ws.Range("E$row").Value = $periodo
    $ws.Range("F$row").Value = $valorMora
    $ws.Range("G$row").Value = $salario
}

Set-Row 16 "CC" "79329212"   "WILLIAM ROJAS SANCHEZ"           "1702" 8000  1000000
Set-Row 17 "CC" "73090118"   "JAIME ALFONSO ARGUELLO NIEBLES"  "1703" 52000 1300000
Set-Row 18 "CC" "1050944522" "ROQUE JOSE SIMANCA VASQUEZ"      "1703" 92000 2500000
Set-Row 19 "CC" "1050944522" "ROQUE JOSE SIMANCA VASQUEZ"      "1703" 8000  2500000
Set-Row 20 "CC" "1047370213" "MARY CARMEN HERNANDEZ MUÑOZ"     "1703" 48000 1200000
Set-Row 21 "CC" "14701504"   "JHON CARLOS TREJOS LOZADA"       "1703" 80000 2000000
Set-Row 22 "CC" "55313486"   "LEIDY KAREN JARAMILLO CASTRO"    "1703" 40000 1000000
Set-Row 23 "CC" "73090118"   "JAIME ALFONSO ARGUELLO NIEBLES"  "1704" 52000 1300000
Set-Row 24 "CC" "1050944522" "ROQUE JOSE SIMANCA VASQUEZ"      "1704" 92000 2500000
Set-Row 25 "CC" "1047370213" "MARY CARMEN HERNANDEZ MUÑOZ"     "1704" 32000 1200000
Set-Row 26 "CC" "14701504"   "JHON CARLOS TREJOS LOZADA"       "1704" 80000 2000000
Set-Row 27 "CC" "55313486"   "LEIDY KAREN JARAMILLO CASTRO"    "1704" 40000 1000000
